$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/15/2025  Through  9/21/2025"

# --- Value-only updates (style/type unchanged) ---
$ws.Range("M14").Value = -85.714285714285
$ws.Range("N14").Value = -90.909090909090
$ws.Range("N15").Value = -83.333333333333
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 91
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = -4.210526315789
$ws.Range("L16").Value = -24.166666666666
$ws.Range("M16").Value = -62.551440329218
$ws.Range("N16").Value = -89.319248826291
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 233.333333333333
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 43.75
$ws.Range("I17").Value = 175
$ws.Range("J17").Value = 225
$ws.Range("K17").Value = -22.222222222222
$ws.Range("L17").Value = -22.222222222222
$ws.Range("M17").Value = -29.149797570850
$ws.Range("N17").Value = -73.484848484848
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 60
$ws.Range("K18").Value = 3.448275862068
$ws.Range("L18").Value = -47.368421052631
$ws.Range("M18").Value = -63.190184049079
$ws.Range("N18").Value = -89.510489510489
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 233.333333333333
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 88.235294117647
$ws.Range("I19").Value = 201
$ws.Range("J19").Value = 204
$ws.Range("K19").Value = -1.470588235294
$ws.Range("L19").Value = -18.623481781376
$ws.Range("M19").Value = -12.608695652173
$ws.Range("N19").Value = -25.278810408921
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 62
$ws.Range("J20").Value = 73
$ws.Range("K20").Value = -15.068493150684
$ws.Range("L20").Value = -23.456790123456
$ws.Range("M20").Value = -7.462686567164
$ws.Range("N20").Value = -86.067415730337
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 108.333333333333
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 58
$ws.Range("H21").Value = 29.310344827586
$ws.Range("I21").Value = 601
$ws.Range("J21").Value = 667
$ws.Range("K21").Value = -9.895052473763
$ws.Range("L21").Value = -25.526641883519
$ws.Range("M21").Value = -38.798370672097
$ws.Range("N21").Value = -79.131944444444
$ws.Range("C23").Value = 3
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 400
$ws.Range("I23").Value = 77
$ws.Range("K23").Value = 13.235294117647
$ws.Range("L23").Value = 16.666666666666
$ws.Range("M23").Value = 28.333333333333
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 53.333333333333
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 65
$ws.Range("H24").Value = 24.615384615384
$ws.Range("I24").Value = 530
$ws.Range("J24").Value = 511
$ws.Range("K24").Value = 3.718199608610
$ws.Range("L24").Value = -16.797488226059
$ws.Range("M24").Value = -7.342657342657
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -42.857142857142
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = -42.105263157894
$ws.Range("I25").Value = 55
$ws.Range("J25").Value = 94
$ws.Range("K25").Value = -41.489361702127
$ws.Range("L25").Value = -65.625
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 77.777777777777
$ws.Range("F26").Value = 39
$ws.Range("H26").Value = 8.333333333333
$ws.Range("I26").Value = 312
$ws.Range("J26").Value = 307
$ws.Range("K26").Value = 1.628664495114
$ws.Range("L26").Value = -13.812154696132
$ws.Range("M26").Value = -47.563025210084
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 500
$ws.Range("J28").Value = 28
$ws.Range("K28").Value = 17.857142857142
$ws.Range("L28").Value = 37.5
$ws.Range("F29").Value = 1
$ws.Range("M29").Value = -71.739130434782
$ws.Range("N29").Value = -90.972222222222
$ws.Range("F30").Value = 1
$ws.Range("M30").Value = -71.794871794871
$ws.Range("N30").Value = -91.40625

# --- Style/type changes (value + underlying text/number type changes) ---
# Each of these needs a number-format paste from a donor cell that already
# carries the target style, because directly assigning a numeric-looking
# string (e.g. "0") or a plain number auto-picks a style that does not
# necessarily match the workbook's existing shared style table.
$ws.Range("G14").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("G14").PasteSpecial(-4122)

$ws.Range("H14").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("H14").PasteSpecial(-4122)

$ws.Range("F15").Value = "'0"
$ws.Range("A15").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("C20").Value = 2
$ws.Range("D20").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("D22").Value = "'0"
$ws.Range("A22").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").Value = "***.*"
$ws.Range("A22").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("D23").Value = "'0"
$ws.Range("A23").Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("E23").Value = "***.*"
$ws.Range("A23").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("F27").Value = "'0"
$ws.Range("A27").Copy()
$ws.Range("F27").PasteSpecial(-4122)

$ws.Range("C28").Value = 1
$ws.Range("F28").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("D28").Value = 1
$ws.Range("F28").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = 0
$ws.Range("H28").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("G29").Value = "'0"
$ws.Range("A29").Copy()
$ws.Range("G29").PasteSpecial(-4122)

$ws.Range("H29").Value = "***.*"
$ws.Range("A29").Copy()
$ws.Range("H29").PasteSpecial(-4122)

$ws.Range("G30").Value = "'0"
$ws.Range("A30").Copy()
$ws.Range("G30").PasteSpecial(-4122)

$ws.Range("H30").Value = "***.*"
$ws.Range("A30").Copy()
$ws.Range("H30").PasteSpecial(-4122)

$excel.CutCopyMode = 0
